# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of old value -> new value for column F on both affected sheets.
$updates = @{
    1276  = 1587
    665   = 669
    361   = 363
    5148  = 5163
    552   = 554
    10029 = 10140
    258   = 262
    554   = 556
    100   = 107
    51    = 70
    761   = 784
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 6)   # Column F
        $current = $cell.Value2
        if ($null -ne $current -and $updates.ContainsKey([int]$current)) {
            $cell.Value2 = $updates[[int]$current]
        }
    }
}
